$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: "Year" -> "Observation_Date"
$ws.Range("A1").Value = "Observation_Date"
